$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.274.05"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.691.23"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'219.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'0.5243"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.38%  "

$ws.Range("D7").Value = "'1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'0.2692"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.44%  "

$ws.Range("D9").Value = "'0.06446"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.24%  "

$ws.Range("D10").Value = "'22.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.38%  "

$ws.Range("D11").Value = "'0.07480"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "

$ws.Range("D12").Value = "1.695.54"
$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "'4.557"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").Value = "'0.5869"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").Value = "'0.000008588"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("D16").Value = "'64.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "26.378.29"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").Value = "'4.982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "'10.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").Value = "'190.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "

$ws.Range("D22").Value = "'6.248"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("D23").Value = "'1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "'145.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.90%  "

$ws.Range("D25").Value = "'7.695"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").Value = "'0.1240"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.69%  "

$ws.Range("D27").Value = "'15.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("D28").Value = "'0.06730"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.31%  "

$ws.Range("D29").Value = "'1.348"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.38%  "

$ws.Range("D30").Value = "'1.333"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("D31").Value = "'3.606"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.45%  "

$ws.Range("D32").Value = "'3.558"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "

$ws.Range("D33").Value = "'1.669"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").Value = "'1.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("D35").Value = "'0.6218"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.30%  "

$ws.Range("D36").Value = "'2.388"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("D37").Value = "'2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.70%  "

$ws.Range("D38").Value = "'6.299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.18%  "

$ws.Range("D39").Value = "'0.01618"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").Value = "1.103.99"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("D41").Value = "'0.8809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.06%  "

$ws.Range("D42").Value = "'1.016"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").Value = "'100.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.53%  "

$ws.Range("D44").Value = "1.842.53"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("D45").Value = "'0.00000000112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").Value = "'56.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.27%  "

$ws.Range("D47").Value = "'8.165"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.46%  "

$ws.Range("D48").Value = "'1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'0.05264"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("D50").Value = "'0.4298"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").Value = "'6.014"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.90%  "
